$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('C2').Value = 'maa://24702 (94.1), maa://25390 (97.42), maa://36681 (92.19)'
$ws.Range('AE2').Value = 'maa://25251 (92.21), ***maa://21730 (17.19), ***maa://39501 (18.18), *maa://36675 (60.0)'
$ws.Range('O3').Value = 'maa://21249 (95.1), maa://26254 (95.24)'
$ws.Range('W3').Value = 'maa://27396 (84.81), maa://27484 (95.65), maa://27480 (82.35)'
$ws.Range('AA3').Value = 'maa://24390 (96.0)'
$ws.Range('S4').Value = 'maa://32509 (98.73), maa://22754 (91.67), maa://27295 (80.39), *maa://21746 (55.81), *maa://31008 (78.05)'
$ws.Range('W4').Value = '**maa://32495 (48.13), ***maa://31785 (16.51), ***maa://36683 (26.67)'
$ws.Range('AE4').Value = '*maa://30062 (63.41), ***maa://26209 (13.04), *maa://39394 (76.92)'
$ws.Range('C5').Value = 'maa://21245 (82.63), maa://22744 (83.33)'
$ws.Range('AA5').Value = '*maa://29863 (74.07), ***maa://22752 (13.33), **maa://26013 (42.86)'
$ws.Range('AE6').Value = '*maa://33152 (58.06), ***maa://22770 (28.57)'
$ws.Range('W7').Value = 'maa://22399 (94.62), *maa://22758 (72.0)'
$ws.Range('C8').Value = '*maa://21476 (69.05), ***maa://39431 (25.0), *maa://37551 (57.14)'
$ws.Range('C10').Value = '***maa://25695 (19.3), **maa://32237 (38.89), ***maa://34206 (14.29), ***maa://39951 (18.18), **maa://39243 (40.0)'
$ws.Range('S10').Value = 'maa://27395 (96.62), maa://22755 (87.5), **maa://22756 (40.91), ***maa://21737 (10.61)'
$ws.Range('W10').Value = 'maa://22301 (97.35), maa://22726 (100.0)'
$ws.Range('AE10').Value = '*maa://25021 (56.34), *maa://22733 (58.62), maa://22761 (100.0)'
$ws.Range('S11').Value = 'maa://22747 (94.33), maa://22501 (98.11)'
$ws.Range('W11').Value = 'maa://36713 (97.81)'
$ws.Range('AA12').Value = 'maa://23669 (95.83), maa://36677 (94.87), maa://39872 (83.33)'
$ws.Range('AE12').Value = '*maa://28932 (78.45), *maa://20106 (63.64), *maa://22769 (62.96)'
$ws.Range('C13').Value = 'maa://24999 (91.37), maa://36673 (91.8), maa://25001 (85.51)'
$ws.Range('G13').Value = '*maa://21248 (75.49), **maa://22728 (47.62)'
$ws.Range('W13').Value = '*maa://34957 (78.57), *maa://22768 (53.33)'
$ws.Range('O14').Value = 'maa://23250 (98.47), maa://20107 (87.1), maa://22772 (100.0), **maa://22745 (50.0)'
$ws.Range('S14').Value = '*maa://22471 (59.42), maa://22521 (94.38)'
$ws.Range('G15').Value = 'maa://24304 (88.4), maa://21478 (91.18)'
$ws.Range('S15').Value = 'maa://23892 (98.63)'
$ws.Range('C16').Value = 'maa://21441 (96.15), maa://36679 (93.55), maa://37650 (95.45)'
$ws.Range('O16').Value = 'maa://28504 (91.49)'
$ws.Range('W16').Value = 'maa://28501 (97.4), maa://28051 (95.83)'
$ws.Range('AA16').Value = 'maa://26228 (96.15)'
$ws.Range('AE16').Value = '*maa://23911 (61.54), maa://27755 (91.67)'
$ws.Range('C17').Value = 'maa://21624 (81.25)'
$ws.Range('G17').Value = 'maa://22430 (88.57), maa://39599 (83.33)'
$ws.Range('C18').Value = 'maa://24570 (96.47)'
$ws.Range('G18').Value = 'maa://24421 (90.48)'
$ws.Range('W18').Value = 'maa://21917 (97.47), maa://22741 (83.33)'
$ws.Range('S19').Value = 'maa://24386 (98.73)'
$ws.Range('AA19').Value = '*maa://30709 (60.48), *maa://36668 (52.17)'
$ws.Range('G20').Value = 'maa://22864 (88.37)'
$ws.Range('K20').Value = 'maa://41331 (89.29)'
$ws.Range('AA21').Value = '*maa://21443 (78.9), **maa://23820 (30.91)'
$ws.Range('W22').Value = 'maa://21282 (98.81), *maa://37649 (64.71)'
$ws.Range('AE22').Value = 'maa://29658 (94.59)'
$ws.Range('K23').Value = 'maa://39756 (91.95), maa://39875 (95.56)'
$ws.Range('W23').Value = '*maa://28503 (62.07)'
$ws.Range('AA23').Value = 'maa://29652 (97.3)'
$ws.Range('C24').Value = 'maa://24368 (80.56)'
$ws.Range('W24').Value = 'maa://23504 (93.08), maa://29988 (85.93), **maa://22892 (40.43), *maa://25141 (76.86), *maa://36663 (79.63), ***maa://22815 (23.08)'
$ws.Range('AE24').Value = 'maa://22523 (84.86), *maa://36672 (76.19), maa://29910 (93.88), **maa://21440 (34.55)'
$ws.Range('G25').Value = '*maa://29063 (76.15), *maa://25311 (73.91), ***maa://22725 (4.84)'
$ws.Range('K25').Value = 'maa://24378 (88.57)'
$ws.Range('AA25').Value = 'maa://31215 (83.75), *maa://24516 (79.07), maa://26001 (88.89)'
$ws.Range('W26').Value = 'maa://24389 (96.0)'
$ws.Range('G27').Value = '**maa://21283 (48.65), maa://34494 (100.0), **maa://36665 (44.44), maa://39601 (87.5)'
$ws.Range('S27').Value = '*maa://30624 (75.68)'
$ws.Range('K28').Value = '*maa://30770 (78.57)'
$ws.Range('W28').Value = 'maa://39929 (85.56), ***maa://39723 (15.15), maa://41749 (81.82)'
$ws.Range('AE28').Value = 'maa://36660 (94.09), *maa://36701 (64.0)'
$ws.Range('K29').Value = 'maa://28432 (93.43), *maa://28440 (72.84), maa://31400 (100.0), *maa://28650 (66.67)'
$ws.Range('AE29').Value = '*maa://24080 (68.78), ***maa://34960 (9.09)'
$ws.Range('S30').Value = '*maa://32940 (66.67), maa://24388 (93.75)'
$ws.Range('W30').Value = '*maa://39477 (71.43)'
$ws.Range('G32').Value = 'maa://21895 (97.01), maa://36667 (98.0), **maa://20793 (38.78), maa://22760 (100.0)'
$ws.Range('K32').Value = 'maa://28065 (94.59)'
$ws.Range('S32').Value = 'maa://41108 (87.5), maa://41238 (94.12)'
$ws.Range('S34').Value = 'maa://24526 (93.07)'
$ws.Range('AE34').Value = '*maa://32650 (69.23)'
$ws.Range('K35').Value = 'maa://41296 (97.3)'
$ws.Range('S35').Value = 'maa://24842 (93.88)'
$ws.Range('AE35').Value = 'maa://39479 (88.89)'
$ws.Range('S36').Value = 'maa://27613 (98.96)'
$ws.Range('O37').Value = 'maa://21280 (89.13), *maa://21239 (72.73)'
$ws.Range('S37').Value = '**maa://39354 (40.0)'
$ws.Range('AE38').Value = 'maa://36697 (85.6)'
$ws.Range('G39').Value = 'maa://25199 (86.11), maa://36670 (89.23), maa://30434 (88.89), ***maa://25036 (16.0)'
$ws.Range('O39').Value = 'maa://24709 (91.92)'
$ws.Range('O41').Value = '**maa://35616 (34.62)'
$ws.Range('S44').Value = 'maa://39366 (90.48)'
$ws.Range('S45').Value = '*maa://39364 (57.14)'
$ws.Range('G46').Value = 'maa://35931 (92.44)'
$ws.Range('O49').Value = '*maa://39643 (71.43)'
$ws.Range('G55').Value = 'maa://32532 (92.27)'
$ws.Range('G57').Value = 'maa://25176 (97.83)'
$ws.Range('G59').Value = 'maa://27746 (84.0), maa://31270 (96.94)'
$ws.Range('G60').Value = '**maa://40438 (42.86)'
